$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so numeric-looking
# strings (e.g. "28.240.89") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.240.89"
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$ws.Range("D3").Value = "1.806.69"
$ws.Range("E3").Value = "  +2.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "338.01"
$ws.Range("E5").Value = "  -0.46%  "

# Row 6
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").Value = "0.4564"
$ws.Range("E7").Value = "  +20.85%  "

# Row 8
$ws.Range("D8").Value = "0.3539"
$ws.Range("E8").Value = "  +4.96%  "

# Row 9
$ws.Range("D9").Value = "45.59"
$ws.Range("E9").Value = "  +0.15%  "

# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.07617"
$ws.Range("E11").Value = "  +4.45%  "

# Row 12
$ws.Range("D12").Value = "22.77"
$ws.Range("E12").Value = "  -1.44%  "

# Row 13
$ws.Range("D13").Value = "0.9991"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("D14").Value = "6.260"
$ws.Range("E14").Value = "  -0.29%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.807.83"
$ws.Range("E15").Value = "  +2.32%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.246"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("D17").Value = "0.00001093"
$ws.Range("E17").Value = "  +3.20%  "

# Row 18
$ws.Range("D18").Value = "0.06679"
$ws.Range("E18").Value = "  +0.84%  "

# Row 19
$ws.Range("D19").Value = "81.92"
$ws.Range("E19").Value = "  +0.74%  "

# Row 20
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.24%  "

# Row 21
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -0.08%  "

# Row 22
$ws.Range("D22").Value = "6.394"
$ws.Range("E22").Value = "  +0.51%  "

# Row 23
$ws.Range("D23").Value = "28.266.14"
$ws.Range("E23").Value = "  +1.56%  "

# Row 24
$ws.Range("D24").Value = "11.99"
$ws.Range("E24").Value = "  +1.33%  "

# Row 25
$ws.Range("D25").Value = "2.402"
$ws.Range("E25").Value = "  +0.55%  "

# Row 26
$ws.Range("D26").Value = "20.78"
$ws.Range("E26").Value = "  +3.47%  "

# Row 27
$ws.Range("D27").Value = "2.427"
$ws.Range("E27").Value = "  +2.77%  "

# Row 28
$ws.Range("D28").Value = "155.47"
$ws.Range("E28").Value = "  +2.14%  "

# Row 29
$ws.Range("D29").Value = "2.012.85"
$ws.Range("E29").Value = "  +2.24%  "

# Row 30
$ws.Range("D30").Value = "1.298"
$ws.Range("E30").Value = "  -13.99%  "

# Row 31
$ws.Range("D31").Value = "133.58"
$ws.Range("E31").Value = "  +0.55%  "

# Row 32
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  +0.71%  "

# Row 33
$ws.Range("D33").Value = "5.958"
$ws.Range("E33").Value = "  +0.24%  "

# Row 34
$ws.Range("D34").Value = "0.09496"
$ws.Range("E34").Value = "  +8.04%  "

# Row 35
$ws.Range("D35").Value = "0.02388"
$ws.Range("E35").Value = "  +0.56%  "

# Row 36
$ws.Range("D36").Value = "12.22"
$ws.Range("E36").Value = "  -2.02%  "

# Row 37
$ws.Range("D37").Value = "0.6750"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2175"
$ws.Range("E38").Value = "  +2.17%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06256"
$ws.Range("E39").Value = "  -0.62%  "

# Row 40
$ws.Range("D40").Value = "5.205"
$ws.Range("E40").Value = "  -0.26%  "

# Row 41
$ws.Range("D41").Value = "1.487"
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
$ws.Range("D42").Value = "1.218"
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("D43").Value = "8.177"
$ws.Range("E43").Value = "  +0.50%  "

# Row 44
$ws.Range("D44").Value = "0.9981"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
$ws.Range("E45").Value = "  +1.19%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.865"
$ws.Range("E46").Value = "  +0.49%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.6125"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
$ws.Range("D48").Value = "129.89"
$ws.Range("E48").Value = "  -1.74%  "

# Row 49
$ws.Range("D49").Value = "2.041"
$ws.Range("E49").Value = "  +0.13%  "

# Row 50
$ws.Range("D50").Value = "0.07113"
$ws.Range("E50").Value = "  -2.26%  "

# Row 51
$ws.Range("D51").Value = "1.165"
$ws.Range("E51").Value = "  -2.07%  "

# Restore default (Normal) style on the Price/Volume columns so the
# workbook styling matches the original (unstyled) cells.
$ws.Range("D2:E51").Style = "Normal"
